$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (column C is never touched by this update, always default style)
$defaultStyle = $ws.Range("C2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.691.54"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  +4.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.373.86"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.61"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +6.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.38"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +3.93%  "

$ws.Range("E9").Value = "  +4.48%  "

$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.29"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +3.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000280"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +7.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "640.72"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  +12.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.911.48"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +1.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.54"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.778.57"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +4.51%  "

$ws.Range("E17").Value = "  +1.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.373.67"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +1.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.96"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.10"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +2.44%  "

$ws.Range("E21").Value = "  +2.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.94"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  +1.63%  "

$ws.Range("E25").Value = "  +3.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.85"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +5.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.76"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +4.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "32.96"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +8.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.68"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +2.79%  "

$ws.Range("E30").Value = "  +1.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "616.70"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +10.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.71"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.988.70"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +6.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.09"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +2.44%  "

$ws.Range("E35").Value = "  +2.68%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.07"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("E38").Value = "  +8.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.31"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +6.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +3.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "33.69"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("E42").Value = "  +2.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.43"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +2.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.343"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +3.18%  "

$ws.Range("E45").Value = "  +3.67%  "

$ws.Range("E46").Value = "  +2.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.60"
$ws.Range("D47").Style = $defaultStyle

$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("E49").Value = "  +9.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.37"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +4.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.81"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +7.37%  "
